$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 7.5
$ws.Range("I2").Value = 1.44
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.33
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8
$ws.Range("AD2").Value = 8.5
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 51
$ws.Range("AG2").Value = 351
$ws.Range("AH2").Value = 7
$ws.Range("AM2").Value = 26
$ws.Range("AN2").Value = 8
$ws.Range("AO2").Value = 34
$ws.Range("AT2").Value = 3.25
$ws.Range("BB2").Value = 126
$ws.Range("G4").Value = 1.95
$ws.Range("I4").Value = 3.9
$ws.Range("J4").Value = 2.5
$ws.Range("L4").Value = 4
$ws.Range("Q4").Value = 1.67
$ws.Range("R4").Value = 2.15
$ws.Range("Z4").Value = 17
$ws.Range("AH4").Value = 15
$ws.Range("AO4").Value = 10
$ws.Range("G5").Value = 3.95
$ws.Range("H5").Value = 2.6
$ws.Range("I5").Value = 2.25
$ws.Range("J5").Value = 4.65
$ws.Range("L5").Value = 3.1
$ws.Range("M5").Value = 1.18
$ws.Range("N5").Value = 4.3
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.47
$ws.Range("W5").Value = 7.3
$ws.Range("X5").Value = 19
$ws.Range("Y5").Value = 15
$ws.Range("AA5").Value = 55
$ws.Range("AB5").Value = 80
$ws.Range("AC5").Value = 4.3
$ws.Range("AH5").Value = 4.65
$ws.Range("AI5").Value = 8.75
$ws.Range("AJ5").Value = 10.75
$ws.Range("AK5").Value = 23
$ws.Range("AL5").Value = 29
$ws.Range("AM5").Value = 65
$ws.Range("AN5").Value = 5.4
$ws.Range("AO5").Value = 26
$ws.Range("AP5").Value = 40
$ws.Range("AQ5").Value = 175
$ws.Range("AR5").Value = 250
$ws.Range("AX5").Value = 13.5
$ws.Range("AY5").Value = 32
$ws.Range("AZ5").Value = 70
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 2.1
$ws.Range("G8").Value = 6.9
$ws.Range("H8").Value = 4.8
$ws.Range("I8").Value = 1.39
$ws.Range("J8").Value = 6
$ws.Range("K8").Value = 2.52
$ws.Range("L8").Value = 1.83
$ws.Range("N8").Value = 10.25
$ws.Range("O8").Value = 1.16
$ws.Range("P8").Value = 4.8
$ws.Range("Q8").Value = 1.5
$ws.Range("R8").Value = 2.52
$ws.Range("S8").Value = 1.28
$ws.Range("T8").Value = 3.45
$ws.Range("U8").Value = 1.7
$ws.Range("V8").Value = 2.05
$ws.Range("W8").Value = 20
$ws.Range("X8").Value = 55
$ws.Range("Y8").Value = 23
$ws.Range("Z8").Value = 175
$ws.Range("AC8").Value = 10.25
$ws.Range("AD8").Value = 10.25
$ws.Range("AE8").Value = 18
$ws.Range("AF8").Value = 70
$ws.Range("AG8").Value = 450
$ws.Range("AH8").Value = 8.25
$ws.Range("AI8").Value = 8.5
$ws.Range("AK8").Value = 10.5
$ws.Range("AL8").Value = 11
$ws.Range("AM8").Value = 23
$ws.Range("AN8").Value = 8.25
$ws.Range("AP8").Value = 32
$ws.Range("AS8").Value = 350
$ws.Range("AT8").Value = 3.45
$ws.Range("AU8").Value = 7.6
$ws.Range("AV8").Value = 55
$ws.Range("AW8").Value = 3.45
$ws.Range("AY8").Value = 14
